$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","I","K","L","N","O")

$data = @{
    2 = @(10.83375004420119; 9.206050745303376; 12.11003470496157; 16.86991607391245; 26.62469591873874; 13.91451702733977; 21.7124971938866; 7.99412584584551; 9.581100416257089; 18.28449048474331; 20.86595182941718)
    3 = @(10.55838159648319; 9.21468708579323; 12.1180957827407; 15.89584955866815; 26.74208508639541; 13.96016024177521; 21.80247108980021; 7.801862018353622; 9.564801563295601; 18.33100484410673; 20.94797316095891)
    4 = @(10.38695315693267; 9.220534309491841; 12.12531601372563; 15.26997757108489; 26.82160149607694; 13.99000810512332; 21.86128307796705; 7.682097198892303; 9.556375263256641; 18.36111388065429; 21.00207030532706)
    5 = @(10.31660522478335; 9.223054449065412; 12.12882995901315; 15.008197319934; 26.85586949816754; 14.00263024293431; 21.88614716625316; 7.632927797011757; 9.553342033030447; 18.37377398241013; 21.02505449487343)
    6 = @(10.30489751778343; 9.223481224585836; 12.12944798808033; 14.96433081551589; 26.86167209690519; 14.00475387145299; 21.89033006804546; 7.624743267682669; 9.552862640349744; 18.37589979020974; 21.02892771937583)
    7 = @(10.38600626361724; 9.220567740276051; 12.12536108873213; 15.26647399323133; 26.82205610563113; 13.99017647292026; 21.86161476806507; 7.681435467497802; 9.556332730543884; 18.36128303718559; 21.00237647576951)
    8 = @(10.73935037072336; 9.208915810978823; 12.11234333229319; 16.5399640634477; 26.66362407549088; 13.92987682574578; 21.74278009231114; 7.928228057813123; 9.575154221442393; 18.30020768122503; 20.89345730935571)
    9 = @(11.40941112799735; 9.190367483778129; 12.10479646684006; 19.00274580682531; 26.41223400538371; 13.82606726248896; 21.53802553639428; 8.395827386399525; 9.624469679273941; 18.19268842142235; 20.70951769365184)
    10 = @(11.88258540307328; 9.179336451824819; 12.1101543420611; 20.67494806633232; 26.26404411891849; 13.75856561105716; 21.40479365760996; 8.725963991539725; 9.668061450355577; 18.12110325351383; 20.59246552999998)
    11 = @(12.09272027312019; 9.174876484465015; 12.11494223109886; 21.3917225636224; 26.20462884300842; 13.72975395357509; 21.3479104398457; 8.872591739734064; 9.689441644198888; 18.09013387841266; 20.54314623936083)
    12 = @(12.171484967488; 9.173267414639955; 12.11709146253032; 21.65686569030329; 26.18328551677913; 13.71911575646049; 21.32690540167527; 8.927556295619597; 9.697756093321367; 18.07863505462246; 20.52503557081846)
    13 = @(12.15455865483793; 9.173610413356965; 12.11661366801707; 21.60004134736742; 26.18783068015524; 13.72139478430653; 21.33140540309634; 8.91574435145702; 9.695955798876101; 18.08110137674845; 20.52891086810927)
    14 = @(12.09921689967195; 9.174742508026927; 12.11511232386106; 21.4136618050453; 26.20284971436777; 13.72887329044873; 21.34617161346517; 8.877125201793076; 9.690121332814078; 18.08918328597073; 20.54164492735346)
    15 = @(12.06521107492214; 9.175446331179584; 12.11423642673252; 21.29868154950795; 26.21220001862088; 13.73348952281408; 21.35528606357142; 8.853395510649193; 9.686575831620363; 18.09416343760683; 20.54951856573994)
    16 = @(11.86874326257704; 9.179639114795769; 12.10988858894745; 20.62722412089977; 26.26808848498786; 13.76048662440592; 21.40858604130027; 8.716305723627823; 9.666695002864788; 18.12315920643881; 20.59576775953806)
    17 = @(11.74685246634207; 9.182353888919318; 12.10782212539117; 20.20408069597325; 26.30442681799713; 13.77753360433218; 21.44223766908421; 8.631259225843408; 9.654892531423767; 18.14135516115; 20.62514669433692)
    18 = @(11.67626668408585; 9.183967912530226; 12.10685492136204; 19.95656407809801; 26.32608004282366; 13.7875169813223; 21.46194378494467; 8.582010987765067; 9.648250258273373; 18.15197117399756; 20.64241439713864)
    19 = @(11.65228789308066; 9.184523433093975; 12.10656551595357; 19.87204792380568; 26.33354048787579; 13.7909278340285; 21.46867615779242; 8.565280970701831; 9.646026546084286; 18.15559139076192; 20.64832442609591)
    20 = @(11.75987791239771; 9.182059460218744; 12.10801920554489; 20.24955283636154; 26.30048062812974; 13.77570046270701; 21.4386191126745; 8.640347274390283; 9.656133826895328; 18.13940263311682; 20.62198098808339)
    21 = @(12.11549461918902; 9.17440782145389; 12.11554419753875; 21.46857628470577; 26.19840684135216; 13.72666928868589; 21.34181989028338; 8.88848415831279; 9.691829172545722; 18.08680323502873; 20.5378892738537)
    22 = @(12.34316665746996; 9.16987210154784; 12.12242060645947; 22.22866616901552; 26.13843607700771; 13.69621072979788; 21.2816766673417; 9.047370700392786; 9.716427839747157; 18.05375863635575; 20.48622678921198)
    23 = @(12.22211049954207; 9.172250489823423; 12.11857201316207; 21.82633154458858; 26.16982494693585; 13.71232202314861; 21.31349075230641; 8.962885806715953; 9.703184496160571; 18.07127351762751; 20.51349821905289)
    24 = @(11.75399068844784; 9.18219240549678; 12.10792941760896; 20.22900810905287; 26.30226232799042; 13.77652865647218; 21.44025394267966; 8.636239666090434; 9.655572191049064; 18.14028488843045; 20.62341102847576)
    25 = @(11.23114801908; 9.194927481968515; 12.10491835507053; 18.34778573295695; 26.47385600954469; 13.8526084386826; 21.59039365689064; 8.271450600796316; 9.609820267497888; 18.2204699314607; 20.75610271232742)
}

foreach ($r in $data.Keys) {
    $rowvals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $rowvals[$i]
    }
}
